$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 192.03334
$ws.Range("I33").Value = 161.31818
$ws.Range("J33").Value = 276.5
$ws.Range("K33").Value = 161.31818
$ws.Range("L33").Value = 276.5
$ws.Range("M33").Value = 67.68181999999999
$ws.Range("N33").Value = -734.5
$ws.Range("H74").Value = 3989.2222
$ws.Range("I74").Value = 3967.6667
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3967.6667
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -3031.6667
$ws.Range("N74").Value = -5872
$ws.Range("H76").Value = 3550.5
$ws.Range("I76").Value = 3410
$ws.Range("K76").Value = 3410
$ws.Range("M76").Value = -3095
$ws.Range("H77").Value = 3989.2222
$ws.Range("I77").Value = 3967.6667
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19838.3335
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -15158.3335
$ws.Range("N77").Value = -29360
$ws.Range("H79").Value = 3550.5
$ws.Range("I79").Value = 3410
$ws.Range("K79").Value = 3410
$ws.Range("M79").Value = -2318
$ws.Range("H129").Value = 874.0405
$ws.Range("J129").Value = 976.87036
$ws.Range("L129").Value = 2930.61108
$ws.Range("N129").Value = -12930.61108
$ws.Range("H135").Value = 54174.105
$ws.Range("I135").Value = 68001.336
$ws.Range("J135").Value = 2322
$ws.Range("K135").Value = 612012.024
$ws.Range("L135").Value = 20898
$ws.Range("M135").Value = -609477.024
$ws.Range("N135").Value = -25968
$ws.Range("H137").Value = 5884679.5
$ws.Range("J137").Value = 2861.111
$ws.Range("L137").Value = 8583.332999999999
$ws.Range("N137").Value = -13683.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2559.875
$ws.Range("I45").Value = 2568.4285
$ws.Range("K45").Value = 2568.4285
$ws.Range("M45").Value = -2191.4285
$ws.Range("H106").Value = 49991.332
$ws.Range("J106").Value = 49991.332
$ws.Range("L106").Value = 49991.332
$ws.Range("N106").Value = -52515.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 36493.715
$ws.Range("J76").Value = 36493.715
$ws.Range("L76").Value = 36493.715
$ws.Range("N76").Value = -37123.715
$ws.Range("H79").Value = 36493.715
$ws.Range("J79").Value = 36493.715
$ws.Range("L79").Value = 36493.715
$ws.Range("N79").Value = -38677.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3262.6553
$ws.Range("I31").Value = 1846
$ws.Range("J31").Value = 4008.2632
$ws.Range("K31").Value = 1846
$ws.Range("L31").Value = 4008.2632
$ws.Range("M31").Value = -1551
$ws.Range("N31").Value = -4598.263199999999
$ws.Range("H34").Value = 3262.6553
$ws.Range("I34").Value = 1846
$ws.Range("J34").Value = 4008.2632
$ws.Range("K34").Value = 1846
$ws.Range("L34").Value = 4008.2632
$ws.Range("M34").Value = -1644
$ws.Range("N34").Value = -4412.263199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3000.625
$ws.Range("I55").Value = 3500
$ws.Range("J55").Value = 2501.25
$ws.Range("K55").Value = 10500
$ws.Range("L55").Value = 7503.75
$ws.Range("M55").Value = -10323
$ws.Range("N55").Value = -7857.75
$ws.Range("H122").Value = 1068.5
$ws.Range("J122").Value = 1215.6666
$ws.Range("L122").Value = 10940.9994
$ws.Range("N122").Value = -15840.9994
$ws.Range("H132").Value = 1154.3549
$ws.Range("I132").Value = 725.5
$ws.Range("J132").Value = 1507.5294
$ws.Range("K132").Value = 6529.5
$ws.Range("L132").Value = 13567.7646
$ws.Range("M132").Value = -3999.5
$ws.Range("N132").Value = -18627.7646
$ws.Range("H133").Value = 2951.6956
$ws.Range("I133").Value = 2478.75
$ws.Range("J133").Value = 4032.7144
$ws.Range("K133").Value = 7436.25
$ws.Range("L133").Value = 12098.1432
$ws.Range("M133").Value = -2376.25
$ws.Range("N133").Value = -22218.1432
$ws.Range("H134").Value = 3699.8572
$ws.Range("J134").Value = 5845.4736
$ws.Range("L134").Value = 17536.4208
$ws.Range("N134").Value = -27676.4208
$ws.Range("H137").Value = 32683.63
$ws.Range("I137").Value = 871.3570999999999
$ws.Range("J137").Value = 121758
$ws.Range("K137").Value = 2614.0713
$ws.Range("L137").Value = 365274
$ws.Range("M137").Value = 2485.9287
$ws.Range("N137").Value = -375474
$ws.Range("H139").Value = 203781.98
$ws.Range("I139").Value = 589618.75
$ws.Range("J139").Value = 5017.5757
$ws.Range("K139").Value = 1768856.25
$ws.Range("L139").Value = 15052.7271
$ws.Range("M139").Value = -1763716.25
$ws.Range("N139").Value = -25332.7271

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 38420
$ws.Range("J100").Value = 38420
$ws.Range("L100").Value = 38420
$ws.Range("N100").Value = -40584
$ws.Range("H101").Value = 63826.668
$ws.Range("J101").Value = 63826.668
$ws.Range("L101").Value = 63826.668
$ws.Range("N101").Value = -70316.66800000001
$ws.Range("H126").Value = 1670.12
$ws.Range("I126").Value = 1354.4117
$ws.Range("J126").Value = 2341
$ws.Range("K126").Value = 4063.2351
$ws.Range("L126").Value = 7023
$ws.Range("M126").Value = -1593.2351
$ws.Range("N126").Value = -11963
$ws.Range("H132").Value = 114268.78
$ws.Range("I132").Value = 113326.22
$ws.Range("J132").Value = 115211.336
$ws.Range("K132").Value = 339978.66
$ws.Range("L132").Value = 345634.008
$ws.Range("M132").Value = -337448.66
$ws.Range("N132").Value = -350694.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2764.7646
$ws.Range("I40").Value = 2500.0625
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 2500.0625
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -2364.0625
$ws.Range("N40").Value = -7272
$ws.Range("H76").Value = 38465.2
$ws.Range("J76").Value = 38465.2
$ws.Range("L76").Value = 38465.2
$ws.Range("N76").Value = -39141.2
$ws.Range("H79").Value = 38465.2
$ws.Range("J79").Value = 38465.2
$ws.Range("L79").Value = 38465.2
$ws.Range("N79").Value = -40805.2
$ws.Range("H94").Value = 31500
$ws.Range("J94").Value = 31500
$ws.Range("L94").Value = 31500
$ws.Range("N94").Value = -32852
$ws.Range("H103").Value = 30457.691
$ws.Range("J103").Value = 30457.691
$ws.Range("L103").Value = 30457.691
$ws.Range("N103").Value = -32801.691
$ws.Range("H132").Value = 104166.2
$ws.Range("I132").Value = 2900
$ws.Range("J132").Value = 147566
$ws.Range("K132").Value = 8700
$ws.Range("L132").Value = 442698
$ws.Range("M132").Value = -6170
$ws.Range("N132").Value = -447758

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 988.25
$ws.Range("I126").Value = 778.84
$ws.Range("J126").Value = 2733.3333
$ws.Range("K126").Value = 2336.52
$ws.Range("L126").Value = 8199.999899999999
$ws.Range("M126").Value = 133.48
$ws.Range("N126").Value = -13139.9999
